$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2147063333333333
$ws.Range("H2").Value = 0.644119
$ws.Range("I2").Value = 0.00381773955517184
$ws.Range("J2").Value = 0.00381773955517184
$ws.Range("M2").Value = 0.1811433333333334
$ws.Range("N2").Value = 0.5434300000000001
$ws.Range("O2").Value = 0.0111261749556462
$ws.Range("P2").Value = 0.01112617495564619
$ws.Range("Q2").Value = 0.03889262090777778
$ws.Range("R2").Value = 0.3500335881700001
$ws.Range("S2").Value = 0.00004247683822593277
$ws.Range("T2").Value = 0.00004247683822593277
$ws.Range("G3").Value = 0.2147063333333333
$ws.Range("H3").Value = 0.644119
$ws.Range("I3").Value = 0.00381773955517184
$ws.Range("J3").Value = 0.00381773955517184
$ws.Range("O3").Value = 0.8246098959508241
$ws.Range("P3").Value = 0.8246098959508241
$ws.Range("Q3").Value = 2.882503664365111
$ws.Range("R3").Value = 25.942532979286
$ws.Range("S3").Value = 0.003148145817357596
$ws.Range("T3").Value = 0.003148145817357597
$ws.Range("G4").Value = 0.2147063333333333
$ws.Range("H4").Value = 0.644119
$ws.Range("I4").Value = 0.00381773955517184
$ws.Range("J4").Value = 0.00381773955517184
$ws.Range("M4").Value = 2.659118666666667
$ws.Range("N4").Value = 7.977356
$ws.Range("O4").Value = 0.1633282272592126
$ws.Range("P4").Value = 0.1633282272592126
$ws.Range("Q4").Value = 0.5709296188182222
$ws.Range("R4").Value = 5.138366569364
$ws.Range("S4").Value = 0.0006235446336835915
$ws.Range("T4").Value = 0.0006235446336835915
$ws.Range("G5").Value = 0.2147063333333333
$ws.Range("H5").Value = 0.644119
$ws.Range("I5").Value = 0.00381773955517184
$ws.Range("J5").Value = 0.00381773955517184
$ws.Range("M5").Value = 0.015234
$ws.Range("N5").Value = 0.045702
$ws.Range("O5").Value = 0.0009357018343171013
$ws.Range("P5").Value = 0.0009357018343171013
$ws.Range("Q5").Value = 0.003270836282
$ws.Range("R5").Value = 0.029437526538
$ws.Range("S5").Value = 0.000003572265904719245
$ws.Range("T5").Value = 0.000003572265904719245
$ws.Range("I6").Value = 0.9856461909412342
$ws.Range("J6").Value = 0.9856461909412343
$ws.Range("M6").Value = 0.1811433333333334
$ws.Range("N6").Value = 0.5434300000000001
$ws.Range("O6").Value = 0.0111261749556462
$ws.Range("P6").Value = 0.01112617495564619
$ws.Range("Q6").Value = 10.04111545575222
$ws.Range("R6").Value = 90.37003910177
$ws.Range("S6").Value = 0.01096647196477843
$ws.Range("T6").Value = 0.01096647196477843
$ws.Range("I7").Value = 0.9856461909412342
$ws.Range("J7").Value = 0.9856461909412343
$ws.Range("O7").Value = 0.8246098959508241
$ws.Range("P7").Value = 0.8246098959508241
$ws.Range("Q7").Value = 744.1913509544628
$ws.Range("R7").Value = 6697.722158590165
$ws.Range("S7").Value = 0.8127736029563772
$ws.Range("T7").Value = 0.8127736029563774
$ws.Range("I8").Value = 0.9856461909412342
$ws.Range("J8").Value = 0.9856461909412343
$ws.Range("M8").Value = 2.659118666666667
$ws.Range("N8").Value = 7.977356
$ws.Range("O8").Value = 0.1633282272592126
$ws.Range("P8").Value = 0.1633282272592126
$ws.Range("Q8").Value = 147.3999459500538
$ws.Range("R8").Value = 1326.599513550484
$ws.Range("S8").Value = 0.1609838450712272
$ws.Range("T8").Value = 0.1609838450712271
$ws.Range("I9").Value = 0.9856461909412342
$ws.Range("J9").Value = 0.9856461909412343
$ws.Range("M9").Value = 0.015234
$ws.Range("N9").Value = 0.045702
$ws.Range("O9").Value = 0.0009357018343171013
$ws.Range("P9").Value = 0.0009357018343171013
$ws.Range("Q9").Value = 0.844449254842
$ws.Range("R9").Value = 7.600043293578
$ws.Range("S9").Value = 0.0009222709488513767
$ws.Range("T9").Value = 0.0009222709488513768
$ws.Range("G10").Value = 0.5925393333333333
$ws.Range("H10").Value = 1.777618
$ws.Range("I10").Value = 0.01053606950359399
$ws.Range("J10").Value = 0.01053606950359399
$ws.Range("M10").Value = 0.1811433333333334
$ws.Range("N10").Value = 0.5434300000000001
$ws.Range("O10").Value = 0.0111261749556462
$ws.Range("P10").Value = 0.01112617495564619
$ws.Range("Q10").Value = 0.1073345499711111
$ws.Range("R10").Value = 0.9660109497400001
$ws.Range("S10").Value = 0.0001172261526418351
$ws.Range("T10").Value = 0.0001172261526418351
$ws.Range("G11").Value = 0.5925393333333333
$ws.Range("H11").Value = 1.777618
$ws.Range("I11").Value = 0.01053606950359399
$ws.Range("J11").Value = 0.01053606950359399
$ws.Range("O11").Value = 0.8246098959508241
$ws.Range("P11").Value = 0.8246098959508241
$ws.Range("Q11").Value = 7.955036878032444
$ws.Range("R11").Value = 71.59533190229199
$ws.Range("S11").Value = 0.008688147177089288
$ws.Range("T11").Value = 0.00868814717708929
$ws.Range("G12").Value = 0.5925393333333333
$ws.Range("H12").Value = 1.777618
$ws.Range("I12").Value = 0.01053606950359399
$ws.Range("J12").Value = 0.01053606950359399
$ws.Range("M12").Value = 2.659118666666667
$ws.Range("N12").Value = 7.977356
$ws.Range("O12").Value = 0.1633282272592126
$ws.Range("P12").Value = 0.1633282272592126
$ws.Range("Q12").Value = 1.575632402000889
$ws.Range("R12").Value = 14.180691618008
$ws.Range("S12").Value = 0.001720837554301858
$ws.Range("T12").Value = 0.001720837554301858
$ws.Range("G13").Value = 0.5925393333333333
$ws.Range("H13").Value = 1.777618
$ws.Range("I13").Value = 0.01053606950359399
$ws.Range("J13").Value = 0.01053606950359399
$ws.Range("M13").Value = 0.015234
$ws.Range("N13").Value = 0.045702
$ws.Range("O13").Value = 0.0009357018343171013
$ws.Range("P13").Value = 0.0009357018343171013
$ws.Range("Q13").Value = 0.009026744203999999
$ws.Range("R13").Value = 0.08124069783599999
$ws.Range("S13").Value = 0.000009858619561005365
$ws.Range("T13").Value = 0.000009858619561005365
